# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
#   - theme2.xml (the theme actually referenced by the slide master / all
#     slides, i.e. the one the user sees and edits via the Design tab) had
#     the "Integral" color scheme (greens/yellows) and ends up with the
#     plain "Office" color scheme (blues/greys).
#   - theme1.xml (only referenced by the Notes Master) had the "Office"
#     color scheme and ends up with the "Integral" one.
# Everything else in both theme parts (font scheme, format scheme) is
# already byte-for-byte identical between the two themes, so the only
# observable difference is the <a:clrScheme> (12 colors) carried by each
# part. This script reproduces that by writing the 12 theme colors of the
# presentation's live theme (ppt/theme/theme2.xml) through the exposed
# PowerPoint object model, switching them from the "Integral" palette to
# the stock "Office" palette - the same end state a user gets after
# picking the built-in "Office" design/color scheme from the Design tab.

$p = $ppt.ActivePresentation

# Any Slide's ThemeColorScheme reaches into the one live DrawingML theme
# that backs the slide master (and therefore every slide) - all 12 slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) are addressable here.
$slide = $p.Slides.Item(1)
$colors = $slide.ThemeColorScheme

# Target palette = the stock Office theme colors (values taken from the
# theme that previously lived in theme1.xml), expressed as the packed
# R + G*256 + B*65536 integers PowerPoint's ColorFormat.RGB expects.
$colors.Item(1).RGB  = 0         # dk1      #000000
$colors.Item(2).RGB  = 16777215  # lt1      #FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      #44546A
$colors.Item(4).RGB  = 15132391  # lt2      #E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  #5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  #ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  #A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  #FFC000
$colors.Item(9).RGB  = 12874308  # accent5  #4472C4
$colors.Item(10).RGB = 4697456   # accent6  #70AD47
$colors.Item(11).RGB = 12673797  # hlink    #0563C1
$colors.Item(12).RGB = 7491477   # folHlink #954F72
